# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps and the "Priority" column
# (sets it to "ht") for the rows that were just handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 11, 13, 14)

foreach ($r in $rows) {
    # Priority column (E) goes from blank to "ht" on both locale sheets.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # Latest Handoff / HO Xliff Generate Date timestamps are refreshed.
    $wsOverview.Range("G$r").Value = "2016-08-23 02:21:54"
    $wsDeDe.Range("H$r").Value     = "2016-08-23 02:21:54"
    $wsZhCn.Range("H$r").Value     = "2016-08-23 02:21:47"
}
